$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H25").Value = 'Everyone is aware of the price of food – from kitchen tables to restaurant counters, food costs impact our lives in profound ways. But another, less understood cost associated with our food is born out of the pressures faced by farmers as they struggle to succeed working in one of the most unforgiving industries in America. Farmers deal with numerous stressors outside of their control including market pressures, weather, price increases, supplies, and many more. Furthermore, farmers deal with the geographic challenge of being isolated from support services in rural areas, and cultural/psychological challenges stemming from personal stigmas against seeking help.<br><br>In this Public Works talk, learn about a unique program in NYS that provides free, confidential support for farmers to help with any number of issues they may face. NY FarmNet approaches farm assistance and support through a “holistic” model by deploying a consultant team representing both the personal and financial aspects of every challenge.<br><br><div class="box"> <img src="images/Oct_4_2023_Adam_Howell_1.jpg" style = ''width: 550'' class=''center''> </div>'
$ws.Range("H26").Value = ' Ghostwriting often gets a bad rap as a cop-out for politicians and celebrities who are too lazy to write their own books. But the truth is, ghostwriting is more common than ever, and for many authors, it''s no longer the shameful secret it once was. So, who actually hires ghostwriters? What do ghostwriters really do? And what''s the future of this ancient profession in a world of generative AI?<br><br><div class="box"> <img src="images/Nov_1_2023_Madison_Fitzpatrick_2.jpg" style = ''width: 550'' class=''center''> </div>'
$ws.Range("H27").Value = ' Ever wonder how others perceive you based on your actions? We all do. Shockingly, recent studies in psychology reveal we''re often quite mistaken about these perceptions. For example, people think they''re less liked by others than they actually are, or they assume others aren''t interested in their unsolicited advice when in fact they are. Such findings could help us get along better with others. But how accurate are these studies? And is it possible that participants in these studies are just trying to appear ''nice'' when answering research questions? In this talk, I''ll share my own research on this topic and discuss if we should take these findings at face value. Join me to discuss the fascinating, messy world of how we''re seen by others.<br><br><div class="box"> <img src="images/Nov_1_2023_Andres_Montealegre_1.jpg" style = ''width: 550'' class=''center''> </div>'

$ws.Range("E26").Select()
